$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting existing rows 4-39 down to 5-40.
$ws.Rows.Item(4).Insert()

# Fill the newly inserted row 4 with the new record's data.
$ws.Cells.Item(4, 1).Value = 8
$ws.Cells.Item(4, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(4, 3).Value = 'Coquimbo'
$ws.Cells.Item(4, 4).Value = (Get-Date -Year 2023 -Month 8 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(4, 4).Style = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
$ws.Cells.Item(4, 5).Value = 4
$ws.Cells.Item(4, 6).Value = 100112013
$ws.Cells.Item(4, 7).Value = 'Alcachofa'
$ws.Cells.Item(4, 8).Value = 'Española'
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 600
$ws.Cells.Item(4, 11).Value = 10000
$ws.Cells.Item(4, 12).Value = 10500
$ws.Cells.Item(4, 13).Value = 10250
$ws.Cells.Item(4, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(4, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(4, 16).Value = 342
$ws.Cells.Item(4, 17).Value = 30
$ws.Cells.Item(4, 18).Value = 'Hortaliza'
